$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(32, 1).Value = 45876
$ws.Cells.Item(32, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(32, 2).Value = 6340
$ws.Cells.Item(32, 3).Value = 6389.7099609375
$ws.Cells.Item(32, 4).Value = 6310.31982421875
$ws.Cells.Item(32, 5).Value = 6374.31982421875
$ws.Cells.Item(32, 6).Value = 5306090000
$ws.Cells.Item(32, 7).Value = -0.0007974800154801
$ws.Cells.Item(32, 8).Value = -1
$ws.Cells.Item(32, 9).Value = 'u tariff on asian mining rig slow bitcoin growth a firm eye canada domestic production outfront medium out q eps drop u auctioning seized m russian superyacht with state room helipad gym spa millennials lead the coffee badging revolt to protest return to office a business push to fill empty seat fortune castleinspired home from metal icon list for million global payment gpn q eps jump blazing new trail minnesota city are considering municipal cannabis store whats going on with ionq share today one gas ogs q net income jump ofsis enforcement overhaul what the july consultation mean for uk sanction compliance titus pressing for information on commodity commission nominee earnings season make it clear ai is the only game in town general dynamic board declares dividend myth about metabolism could be holding you back allegiant algt q earnings call transcript sm entertainment revenue up on new release from riize nct wish strong merch sale chatham cldt eps drop wall street oil gas lending down ytd niit learning system limited niit mt announces q fy apriljune result why opendoor technology stock wa tumbling today rosen a ranked and leading law firm encourages irobot corporation investor to secure counsel before important deadline in security class action irbt deloitte india and embark form a strategic alliance to offer endtoend solution for global capability centre how workplace toxicity played a role in the titan submersible disaster mcdonalds announces nostalgic mcdonaldland meal organic avocado mattress are up to off ahead of labor day phone jewelry linen which product could cost more due to trump india tariff jeff wood join patriot growth insurance service a premier client practice leader claires known for piercing million of teen ear file for chapter nd time since a catastrophic funding model court sends sec back to the drawing board tyson beck celebrity mint release first ever coin design featuring basketball icon dennis rodman new york city board denies mayor adam m in matching campaign fund best hp business laptop in strong pick for smart work beloved brooklyn restaurant from chef baxtrom close august claires known for piercing million of teen ear again file for chapter cliq digital ag dylan medium to vote against proposed share buyback and cliq digital currently no longer considers delisting latest global arabinogalactan market sizeshare worth usd million by at a cagr custom market insight analysis outlook leader report trend forecast segmentation growth rate value swot analysis coupa add tariff impact planning to supply chain tool amgen stock sink a market eye maritides next move how can you borrow less for grad school live will q earnings fuel doordashs dash next rally jenavalve responds to ftc action against proposed edward lifesciences acquisition roger williams university school of law named to best law school for by the princeton review israeli society is living in trauma say frenchisraeli journalist pioneer realty capital secures million dollar bond refinancing for texas student housing project from event to improvement step for turning dash cam alert into driver coaching a cautionary tale state ag prevail with a lump of coal for major investment firm a rally for apple lead wall street higher tween jewelry retailer claires file for second bankruptcy bloombergcom brazil request consultation at world trade organization over trump tariff brazil request consultation at world trade organization over trump tariff smart payment solution launch mony a revolutionary upi app for tourist and nris visiting india apple share pop ahead of trumpcook investment announcement at white house best crypto to buy today arctic pablos entry might be the last golden ticket a shiba inu and cheems revive interest c o r r e c t i o n enercom inc enercom announces andrew rapp senior advisor in the u department of energy a keynote speaker on august th at enercoms th anniversary energy investment conference rhode island writes it own cybersecurity rule for nonbanks ap trending summarybrief at pm edt rosen skilled investor counsel encourages lockheed martin corporation investor to secure rosen skilled investor counsel encourages lockheed martin corporation investor to secure counsel before important deadline in security class action lmt china fight mosquitoborne chikungunya virus with drone fine and net a thousand fall ill a preview of sabre earnings preview relay therapeutics earnings direct equity source featured in exclusive qa with two time emmy award winner jim knox a look ahead realreals earnings forecast claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since claires known for piercing million of teen ear file for chapter nd time since a look ahead advanced drainage systemss earnings forecast earnings outlook for warby parker snap stock is falling wednesday whats going on u is auctioning a seized m russian superyacht with state room a helipad a gym and a spa cloudsquare unveils new website and brand direction cementing it role a the leading lending platform for alternative finance we want to get there first and claim that for america nasa chief explains push for nuclear reactor on the moon video space walker dunlops earnings a preview exploring construction partnerss earnings expectation tesla cfo taneja excfo kirkhorn named in robotaxi shareholder suit altura join kuviai to enable ai agent in blockchain gaming earnings preview for peloton interactive mapped america hottest zip code revealed region by region realtorcom a peek at dwave quantum future earnings preview penn entertainment earnings rdo equipment co announces layoff delta exchange launch tracker a gamechanging alternative to spot crypto trading behind the scene of caterpillar latest option trend volatility isnt just a metric anymore it driving market structure cava group option frenzy what you need to know jen pawol will be major league baseball first woman to umpire in a regularseason game this weekend mother issue heartbreaking warning after son dy taking popular natural supplement ladbible best hacking surveillance game in origis energy and pioneer community energy announce year solar battery storage agreement richmond american community hosting open house in saratoga spring glencore scrap london listing move missed binance coin bnb roi analyst say ruvi ai ruvi could be your last chance to catch similar gain a phase nears end ethena gain vechain expands but cold wallet cwt token could be the hidden x crypto for adam zyglis cook getting messy suspected south american gang member arrested in maine by border patrol opinion with trump and ice involved los angeles shouldnt host the olympics midair skydiving collision over flathead county leaf one dead moran and marshall lied by omission myth about metabolism could be holding you back pekin restaurant currently closed due to unexpected plumbing issue should professional athlete still be role model the ethical life podcast return national park service to reinstall confederate statue felled by blm mob in should professional athlete still be role model the ethical life podcast gop rep cory mill accused by exgirlfriend of threatening to release sexually explicit image video of her cnn politics texas democrat loses touch with reality compare redistricting to the holocaust rosen a ranked and leading law firm encourages irobot corporation investor to secure counsel before important deadline in security class action irbt cincinnati brawl victim want onlooker who didnt help prosecuted rfk jr launch crackdown on kratom legal morphine substance often found in smoke shop product jonah goldberg kamala harris wont cure what ail the democratic party tim burchett signal hed take blackburn senate seat if she becomes governor once again demanding answer mahmoud khalils terrorist support and his unjustified presence in america former state worker and environmental advocate launch independent run for maine governor roger williams university school of law named to best law school for by the princeton review israeli society is living in trauma say frenchisraeli journalist ucla in talk with trump administration to end m freeze analysis lebanon decision on weapon corner hezbollah watch fox news rave about trump weird white house roof walk james city county police announces school zone speed camera installation former ct criminal justice official pay to resolve ethic complaint denies wrongdoing rosen skilled investor counsel encourages lockheed martin corporation investor to secure counsel before important deadline in security class action lmt china fight mosquitoborne chikungunya virus with drone fine and net a thousand fall ill the latest putin host trump envoy witkoff for crucial talk on ukraine peace deal rep jolanda jones apologizes for comparing texas redistricting plan to holocaust political analyst say if kamala slam biden in her book his team will level her crisp trump suit perfectly our national cognitive dissonance trump say he may deploy national guard in dc editorial redistricting a cynical power play by politician rowley six month into his presidency donald trump ha created a police state se cupp to trump truth is only what he want it be trump envoy witkoff meet putin ahead of russiaukraine peace deadline trump envoy witkoff meet putin ahead of russiaukraine peace deadline jen pawol will be major league baseball first woman to umpire in a regularseason game this weekend the death of highspeed rail not so fast editorial ai help rapidly crack yearlong mystery in alp platte county senator nominated by trump to serve in administration ocean city to consider declaring condemned condo in need of rehabilitation what are cool roof and how do they work india paid the price for underestimating china pl missile disgusting tom homan skewer farleft dem declaring her allegiance to foreign nation over her own country will doj crackdown prompt sanctuary jurisdiction to alter immigration policy democratic politician at iowa state fair john seiler next state school chief will follow delinquent thurmond despite court order federal official use penske truck to trick and arrest immigrant in la several state threaten to redraw congressional map after texas kick off fight several state threaten to redraw congressional map after texas kick off fight several state threaten to redraw congressional map after texas kick off fight several state threaten to redraw congressional map after texas kick off fight several state threaten to redraw congressional map after texas kick off fight several state threaten to redraw congressional map after texas kick off fight several state threaten to redraw congressional map after texas kick off fight democracy be damned texas and california plot dueling congressional gerrymander raap policy study release solving for access and affordability pdabs are not the answer several state threaten to redraw congressional map after texas kick off fight soldier shot before arrest made at army fort stewart in georgia soldier wounded in mass shooting at fort stewart in georgia independentminor party candidate join luzerne county ballot to woody allen epstein wa dracula serviced by young female vampire u sanction narcorapper el makabelico who laundered streaming proceeds to bloodthirsty mexican drug lord a massive u measles outbreak ha slowed but the start of the school year brings renewed risk of spread stanford student newspaper sue trump administration over use of immigration law to target propalestinian student why moderate dems scapegoat bibi dems should drop green energy and other commentary appeal court directs jeffco judge to reconsider releasing man car from lakewood police custody how trump is exploiting musk big ball pritzker democrat vow to win redistricting battle in the wild of the anaconda range expert say elusive gunman ha advantage over searcher pritzker democrat vow to win redistricting battle copy missouri gop see redistricting call a opening for limiting constitutional amendment rising dem star say party isnt playing dead amid gerrymandering chaos from gop bully fort stewart lockdown five soldier shot in active shooter incident dozen killed seeking aid in gaza a israel weighs further military action lawrence odonnell shred trump for playing dumb on ghislaine maxwell prison transfer drill sergeant are back hegseth reinstates practice to make basic great again bridge result for aug tribal group assert sovereignty a fed crack down on genderaffirming care tribal group assert sovereignty a fed crack down on genderaffirming care what to know about mrna vaccine what to know about mrna vaccine lunar lunacy cbs morning host worried colonization could harm the moon indigenous population wut openai is practically giving chatgpt to the government for free trump hate all democrat not just mamdani reader commentary leftist agitator disrupt trump official at student event in washington dc the mindful minute by jon heydenreich funny dozen killed seeking aid in gaza a israel weighs further military action dozen killed seeking aid in gaza a israel weighs further military action indycar distance itself from speedway slammer immigration detention facility defiance over compliance harvard med school continuing to use racial preference complaint say opinion trump destroys our source of information about job this is beyond irresponsible how form of pay got issued to exst louis comptroller mistake and bigger mistake the u dragon lady celebrated it th anniversary with a historic flight america is prepared for a strike on our nuclear weapon deadline alert faruqi faruqi llp investigates claim on behalf of investor of altimmune the real reason democrat are panicking about redistricting support at work central ohio workforce development network launch employee resource network to help employee thrive'
